# Daily attendance processing - 2025-10-26 12:34:54
# Normalizes the "Recorded By" column (G): entries equal to "System"/"system"
# are moved to the front of the comma-separated list, preserving the
# relative order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }
    if ($val -notlike "*,*") {
        continue
    }

    $parts = $val -split ','
    $sysParts = @()
    $otherParts = @()

    foreach ($p in $parts) {
        $trimmed = $p.Trim()
        if ($trimmed.ToLower() -eq "system") {
            $sysParts += $trimmed
        } else {
            $otherParts += $trimmed
        }
    }

    if ($sysParts.Count -eq 0) {
        continue
    }

    $newParts = $sysParts + $otherParts
    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
